$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 1407222.1
$ws.Range("I64").Value = 2562055.2
$ws.Range("J64").Value = 4924.643
$ws.Range("K64").Value = 2562055.2
$ws.Range("L64").Value = 4924.643
$ws.Range("M64").Value = -2561807.2
$ws.Range("N64").Value = -5420.643
$ws.Range("H67").Value = 1407222.1
$ws.Range("I67").Value = 2562055.2
$ws.Range("J67").Value = 4924.643
$ws.Range("K67").Value = 2562055.2
$ws.Range("L67").Value = 4924.643
$ws.Range("M67").Value = -2561197.2
$ws.Range("N67").Value = -6640.643
$ws.Range("H100").Value = 1537.5714
$ws.Range("I100").Value = 1355.3334
$ws.Range("J100").Value = 1674.25
$ws.Range("K100").Value = 1355.3334
$ws.Range("L100").Value = 1674.25
$ws.Range("M100").Value = -814.3334
$ws.Range("N100").Value = -2756.25

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4847.421
$ws.Range("I32").Value = 2283.276
$ws.Range("K32").Value = 2283.276
$ws.Range("M32").Value = -1996.276
$ws.Range("H46").Value = 30991.25
$ws.Range("I46").Value = 52482.5
$ws.Range("J46").Value = 9500
$ws.Range("K46").Value = 52482.5
$ws.Range("L46").Value = 9500
$ws.Range("M46").Value = -52163.5
$ws.Range("N46").Value = -10138
$ws.Range("H61").Value = 11958.921
$ws.Range("I61").Value = 8934.846
$ws.Range("K61").Value = 8934.846
$ws.Range("M61").Value = -8722.846
$ws.Range("H74").Value = 1122.5769
$ws.Range("I74").Value = 572.1429000000001
$ws.Range("J74").Value = 1764.75
$ws.Range("K74").Value = 572.1429000000001
$ws.Range("L74").Value = 1764.75
$ws.Range("M74").Value = 301.8570999999999
$ws.Range("N74").Value = -3512.75
$ws.Range("H77").Value = 1122.5769
$ws.Range("I77").Value = 572.1429000000001
$ws.Range("J77").Value = 1764.75
$ws.Range("K77").Value = 2860.7145
$ws.Range("L77").Value = 8823.75
$ws.Range("M77").Value = 1507.2855
$ws.Range("N77").Value = -17559.75
$ws.Range("H132").Value = 4931.2285
$ws.Range("I132").Value = 1980.5769
$ws.Range("K132").Value = 5941.7307
$ws.Range("M132").Value = -3411.7307
$ws.Range("H136").Value = 11958.921
$ws.Range("I136").Value = 8934.846
$ws.Range("K136").Value = 26804.538
$ws.Range("M136").Value = -24254.538
$ws.Range("H140").Value = 111082.164
$ws.Range("J140").Value = 111082.164
$ws.Range("L140").Value = 111082.164
$ws.Range("N140").Value = -121442.164

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H62").Value = 299999
$ws.Range("J62").Value = 299999
$ws.Range("L62").Value = 299999
$ws.Range("N62").Value = -301371
$ws.Range("H63").Value = 64999
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H65").Value = 299999
$ws.Range("J65").Value = 299999
$ws.Range("L65").Value = 899997
$ws.Range("N65").Value = -906861
$ws.Range("H66").Value = 64999
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 60000
$ws.Range("M66").ClearContents()
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384
$ws.Range("H107").Value = 1772
$ws.Range("I107").Value = 2873.8572
$ws.Range("K107").Value = 2873.8572
$ws.Range("M107").Value = -953.8571999999999
$ws.Range("H134").Value = 3062.862
$ws.Range("I134").Value = 2596.875
$ws.Range("K134").Value = 7790.625
$ws.Range("M134").Value = -5255.625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 378
$ws.Range("J22").Value = 343.4
$ws.Range("L22").Value = 343.4
$ws.Range("N22").Value = -1043.4
$ws.Range("H31").Value = 2201.3142
$ws.Range("I31").Value = 919.86664
$ws.Range("K31").Value = 919.86664
$ws.Range("M31").Value = -624.86664
$ws.Range("H34").Value = 2201.3142
$ws.Range("I34").Value = 919.86664
$ws.Range("K34").Value = 919.86664
$ws.Range("M34").Value = -717.86664
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("N63").Value = 0
$ws.Range("H64").Value = 18999.334
$ws.Range("J64").Value = 18999.334
$ws.Range("L64").Value = 18999.334
$ws.Range("N64").Value = -19495.334
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("N66").Value = 0
$ws.Range("H67").Value = 18999.334
$ws.Range("J67").Value = 18999.334
$ws.Range("L67").Value = 18999.334
$ws.Range("N67").Value = -20715.334
$ws.Range("H108").Value = 76247.25
$ws.Range("J108").Value = 76247.25
$ws.Range("L108").Value = 76247.25
$ws.Range("N108").Value = -83927.25
$ws.Range("H134").Value = 2147.8044
$ws.Range("I134").Value = 2180.4878
$ws.Range("K134").Value = 6541.4634
$ws.Range("M134").Value = -4006.4634

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 4333.25
$ws.Range("I82").Value = 2666.5
$ws.Range("J82").Value = 6000
$ws.Range("K82").Value = 7999.5
$ws.Range("L82").Value = 18000
$ws.Range("M82").Value = -7593.5
$ws.Range("N82").Value = -18812
$ws.Range("H85").Value = 4333.25
$ws.Range("I85").Value = 2666.5
$ws.Range("J85").Value = 6000
$ws.Range("K85").Value = 7999.5
$ws.Range("L85").Value = 18000
$ws.Range("M85").Value = -6595.5
$ws.Range("N85").Value = -20808
$ws.Range("H94").Value = 9678.429
$ws.Range("J94").Value = 9969.833000000001
$ws.Range("L94").Value = 29909.499
$ws.Range("N94").Value = -31261.499
$ws.Range("H96").Value = 12250
$ws.Range("J96").Value = 13000
$ws.Range("L96").Value = 39000
$ws.Range("N96").Value = -43118
$ws.Range("H103").Value = 10792.143
$ws.Range("I103").Value = 136.25
$ws.Range("J103").Value = 25000
$ws.Range("K103").Value = 408.75
$ws.Range("L103").Value = 75000
$ws.Range("M103").Value = 470.25
$ws.Range("N103").Value = -76758
$ws.Range("J131").Value = 6805597.5
$ws.Range("L131").Value = 20416792.5
$ws.Range("N131").Value = -20426872.5
$ws.Range("H138").Value = 58920.277
$ws.Range("I138").Value = 102923
$ws.Range("K138").Value = 308769
$ws.Range("M138").Value = -303629

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1085.4445
$ws.Range("I97").Value = 1180.625
$ws.Range("J97").Value = 324
$ws.Range("K97").Value = 1180.625
$ws.Range("L97").Value = 324
$ws.Range("M97").Value = -684.625
$ws.Range("N97").Value = -1316
$ws.Range("H102").Value = 11774.857
$ws.Range("I102").Value = 12804.454
$ws.Range("K102").Value = 12804.454
$ws.Range("M102").Value = -11182.454

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5422.525
$ws.Range("I7").Value = 3959.5386
$ws.Range("K7").Value = 3959.5386
$ws.Range("M7").Value = -3847.5386
$ws.Range("H93").Value = 1416.6666
$ws.Range("I93").Value = 1250
$ws.Range("J93").Value = 1750
$ws.Range("K93").Value = 1250
$ws.Range("L93").Value = 1750
$ws.Range("M93").Value = -2
$ws.Range("N93").Value = -4246
$ws.Range("H126").Value = 5422.525
$ws.Range("I126").Value = 3959.5386
$ws.Range("K126").Value = 11878.6158
$ws.Range("M126").Value = -9408.6158

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 17980.5
$ws.Range("I74").Value = 17973
$ws.Range("K74").Value = 17973
$ws.Range("M74").Value = -17037
$ws.Range("H77").Value = 17980.5
$ws.Range("I77").Value = 17973
$ws.Range("K77").Value = 53919
$ws.Range("M77").Value = -49239
